$wb = $excel.ActiveWorkbook

# Update "展览" sheet (F2: 73 -> 75, F3: 4 -> 5)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 75
$ws1.Range("F3").Value = 5

# Update "全部类型" sheet (F2: 73 -> 75, F3: 4 -> 5)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 75
$ws4.Range("F3").Value = 5
